# Apply the LinuxForHealth re-branding edit to the immigration-status
# StructureDefinition workbook.
#
# Sheet "Metadata" (first sheet): update URL, Version, Date and Publisher
# values.
# Sheet "Elements" (second sheet): clear the "Constraint(s)" value that was
# recorded against the root "Extension" element row (it now only appears
# against the "Extension.extension" row).

$wb = $excel.ActiveWorkbook

$wsMetadata = $wb.Worksheets.Item("Metadata")
$wsMetadata.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/immigration-status"
$wsMetadata.Range("B3").Value = "8.0.0"
$wsMetadata.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$wsMetadata.Range("B9").Value = "LinuxForHealth Team"

$wsElements = $wb.Worksheets.Item("Elements")
$wsElements.Range("AI2").Value = ""
